# Generate Report for Handoff
# Replace the old source-file GUID/hash identifiers with the newly generated
# ones, and bump the handoff timestamps, across the Overview/zh-cn/de-de
# sheets. Hyperlink targets (Address) are left untouched -- only the
# display text shown in the cell / hyperlink changes.

$wb = $excel.ActiveWorkbook

$oldGuid = "c1b5d003-50e4-4c2a-aa21-532f544a9eb1"
$newGuid = "dfd096ee-948e-4db1-9687-dda1ece87f01"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

$oldZhXlf = "$oldGuid.1e163b1edb66c07a6c001754b531ea12ba8f71b5.zh-cn.xlf"
$newZhXlf = "$newGuid.38819295ffc34122546952c586aec66288bd7449.zh-cn.xlf"

$oldDeXlf = "$oldGuid.1e163b1edb66c07a6c001754b531ea12ba8f71b5.de-de.xlf"
$newDeXlf = "$newGuid.38819295ffc34122546952c586aec66288bd7449.de-de.xlf"

$oldZhDatetime = "2016-03-07 02:36:26"
$newZhDatetime = "2016-03-07 02:37:10"

$oldDeDatetime = "2016-03-07 02:36:37"
$newDeDatetime = "2016-03-07 02:37:20"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/0b31996e71e0dc3a702b26b5d095e17e9773bed2/e2e/$oldMdName"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/0b31996e71e0dc3a702b26b5d095e17e9773bed2/.localization-config"
$configDisplay = ".localization-config"

$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa63631765ddc0e7ac2e9f195afaacd02a4fa152/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$oldZhXlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4fc3aa032273a1fa943050d0db82f8223b68f9dd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$oldDeXlf"

# ---------------------------------------------------------------------------
# Overview sheet: A2 (source file link), A3 (.localization-config link)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMdName

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddress, "", "", $configDisplay)

# ---------------------------------------------------------------------------
# zh-cn sheet: A2 (source file link), C2 (handoff xlf link + datetime)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("C2").Value = $newZhXlf
$wsZhCn.Range("D2").Value = $newZhDatetime

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhXlfAddress, "", "", $newZhXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configAddress, "", "", $configDisplay)

# ---------------------------------------------------------------------------
# de-de sheet: A2 (source file link), C2 (handoff xlf link + datetime)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("C2").Value = $newDeXlf
$wsDeDe.Range("D2").Value = $newDeDatetime

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deXlfAddress, "", "", $newDeXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configAddress, "", "", $configDisplay)
